$d = $word.ActiveDocument

# --- Paragraph 1 (Heading1): "Welcome to the Jungle!!!" -> "Des Kaisers Kleider" ---
# Replace the whole heading text in one shot so Word's run-merge logic also
# cleans up the spell-check <w:proofErr> markers that used to straddle the
# "to"/"the"/"Jungle" runs.
$d.Content.Find.Execute("Welcome to the Jungle!!!", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Des Kaisers Kleider", 2)

# Re-split the merged run back into four runs - "Des K" / "a" / "isers " / "Kleider" -
# matching the run boundaries shown in the target markup. Toggling a character
# property on and back off is a reliable way to force a run split without
# leaving any formatting residue behind.
$p1Start = $d.Paragraphs(1).Range.Start

$part1 = $d.Range($p1Start, $p1Start + 5)      # "Des K"
$part1.Font.Bold = $true
$part1.Font.Bold = $false

$part2 = $d.Range($p1Start + 5, $p1Start + 6)  # "a"
$part2.Font.Bold = $true
$part2.Font.Bold = $false

$part3 = $d.Range($p1Start + 6, $p1Start + 12) # "isers "
$part3.Font.Bold = $true
$part3.Font.Bold = $false

# --- Remove the blank paragraph that used to follow the heading ---
$d.Paragraphs(2).Range.Delete()

# --- Remove the "Es kann nur einen geben!!!" text, keeping the bookmark ---
$d.Content.Find.Execute("Es kann nur einen geben!!!", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "", 2)

# Merge the now-empty paragraph (which still holds the _GoBack bookmark) back
# into the heading paragraph by deleting the heading's end-of-paragraph mark.
$p1End = $d.Paragraphs(1).Range.End
$d.Range($p1End - 1, $p1End).Delete()

# Deleting that mark makes the merged paragraph inherit the second
# paragraph's (Normal) formatting, so restore the Heading1 style explicitly.
$d.Paragraphs(1).Style = "Heading1"
